# Apply the edit described by the diff:
#  - Delete the (empty) column L, shifting columns M:Q left to L:P
#  - Change the active selection on Sheet1 to A28
#  - Set the worksheet to print in Portrait orientation
#  - Turn off "Recalculate before Save" (calcOnSave) so the workbook
#    calc settings match the target (calcPr without calcOnSave="0")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the empty column L -- everything to its right (M:Q) shifts left.
$ws.Columns("L").Delete()

# Update the selected cell/range to match the target file.
$ws.Range("A28").Select()

# Set page orientation to Portrait (adds <pageSetup orientation="portrait".../>).
$ws.PageSetup.Orientation = 1

# Match calcPr change: disable "recalculate before save".
$wb.CalculateBeforeSave = $false
